$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.538.18"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "'1.889.51"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'244.32"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.4777"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("D8").Value = "'0.2895"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06548"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'21.36"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07778"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.886.83"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7366"
$ws.Range("E13").Value = "  +6.63%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'96.26"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.172"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "'277.59"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'30.559.19"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.60"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007620"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'0.9994"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.148.93"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.307"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.214"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.307"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'165.20"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.996"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.381"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.09974"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.513"
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.338"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.105"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04777"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7045"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01860"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.770"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.513"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'70.82"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.928"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4197"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8455"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9994"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'102.92"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.449"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.149"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'929.65"
$ws.Range("E49").Value = "  -4.81%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'35.39"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3876"
$ws.Range("E51").Value = "  +1.32%  "
